$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text would otherwise be auto-recognized as a plain
# number by Excel (single decimal point, e.g. "225.98"). Force them to
# remain plain text (matching the original inlineStr cells) by temporarily
# applying a Text number format, then clearing formatting again so no
# style index is left behind on the cell.
$textForcedCells = @(
    "D5",
    "D11",
    "D15",
    "D18",
    "D21",
    "D25",
    "D27",
    "D32",
    "D33",
    "D40",
    "D44",
    "D50"
)
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D5").Value = "225.98"
$ws.Range("D11").Value = "0.0946"
$ws.Range("D15").Value = "0.624"
$ws.Range("D18").Value = "68.07"
$ws.Range("D21").Value = "10.99"
$ws.Range("D25").Value = "161.69"
$ws.Range("D27").Value = "16.33"
$ws.Range("D32").Value = "3.72"
$ws.Range("D33").Value = "3.65"
$ws.Range("D40").Value = "80.19"
$ws.Range("D44").Value = "13.51"
$ws.Range("D50").Value = "106.04"

foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).ClearFormats()
}

# --- Remaining cells: new text is not number-like (multiple dots, percent
# strings with signs/spaces, or subscript digits), so a plain Value
# assignment keeps them as text without Excel reinterpreting them.
$ws.Range("D2").Value = "34.202.05"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "1.784.71"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "2.041.98"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "1.782.96"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").Value = "34.114.46"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("E19").Value = "  +3.57%  "
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  +4.59%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("E26").Value = "  +2.88%  "
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("E32").Value = "  +3.74%  "
$ws.Range("E33").Value = "  +4.83%  "
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").Value = "1.445.55"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("E37").Value = "  +9.86%  "
$ws.Range("E38").Value = "  +4.08%  "
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("E45").Value = "  +4.54%  "
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "0.0₆0136"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").Value = "1.944.17"
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  +0.10%  "
